$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View state: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 175
[void]$ws.Range("F14").Select()

# --- Column D width ---
$ws.Columns.Item(4).ColumnWidth = 9.5

# --- Header row: add column D "s per epoch" ---
$ws.Range("D1").Value = "s per epoch"

# --- Rename shared text RIALLENA -> RISCRIVI on K4, clear K2 ---
$ws.Range("K2").ClearContents()

# --- Update existing rows 2-6 and add new rows 7-11 ---
# Row 2
$ws.Range("A2").Value = "11n"
$ws.Range("B2").Value = 128
$ws.Range("C2").Value = "no"
$ws.Range("D2").Formula = "=(F2*3600)/E2"
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 0.321
$ws.Range("G2").Value = 0.972
$ws.Range("H2").Value = 0.958
$ws.Range("I2").Value = 0.982
$ws.Range("J2").Value = 0.829

# Row 3
$ws.Range("A3").Value = "11n"
$ws.Range("B3").Value = 256
$ws.Range("C3").Value = "no"
$ws.Range("D3").Formula = "=(F3*3600)/E3"
$ws.Range("E3").Value = 89
$ws.Range("F3").Value = 0.299
$ws.Range("G3").Value = 0.991
$ws.Range("H3").Value = 0.987
$ws.Range("I3").Value = 0.994
$ws.Range("J3").Value = 0.918

# Row 4
$ws.Range("A4").Value = "11n"
$ws.Range("B4").Value = 512
$ws.Range("C4").Value = "no"
$ws.Range("D4").Formula = "=(F4*3600)/E4"
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 0.173
$ws.Range("G4").Value = 0.984
$ws.Range("H4").Value = 0.987
$ws.Range("I4").Value = 0.993
$ws.Range("J4").Value = 0.922
$ws.Range("K4").Value = "RISCRIVI"

# Row 5
$ws.Range("A5").Value = "11n"
$ws.Range("B5").Value = 640
$ws.Range("C5").Value = "no"
$ws.Range("D5").Formula = "=(F5*3600)/E5"
$ws.Range("E5").Value = 71
$ws.Range("F5").Value = 0.541
$ws.Range("G5").Value = 0.997
$ws.Range("H5").Value = 0.992
$ws.Range("I5").Value = 0.995
$ws.Range("J5").Value = 0.954

# Row 6
$ws.Range("A6").Value = "11n"
$ws.Range("B6").Value = 1024
$ws.Range("C6").Value = "no"
$ws.Range("D6").Formula = "=(F6*3600)/E6"
$ws.Range("E6").Value = 78
$ws.Range("F6").Value = 1.03
$ws.Range("G6").Value = 0.994
$ws.Range("H6").Value = 0.997
$ws.Range("I6").Value = 0.995
$ws.Range("J6").Value = 0.962

# Row 7
$ws.Range("A7").Value = "11s"
$ws.Range("B7").Value = 640
$ws.Range("C7").Value = "no"
$ws.Range("D7").Formula = "=(F7*3600)/E7"
$ws.Range("E7").Value = 40
$ws.Range("F7").Value = 0.311
$ws.Range("G7").Value = 0.994
$ws.Range("H7").Value = 0.991
$ws.Range("I7").Value = 0.994
$ws.Range("J7").Value = 0.952

# Row 8
$ws.Range("A8").Value = "11n"
$ws.Range("B8").Value = 640
$ws.Range("C8").Value = "light"
$ws.Range("D8").Formula = "=(F8*3600)/E8"
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = 0.713
$ws.Range("G8").Value = 0.997
$ws.Range("H8").Value = 0.998
$ws.Range("I8").Value = 0.995
$ws.Range("J8").Value = 0.963

# Row 9
$ws.Range("A9").Value = "11n"
$ws.Range("B9").Value = 640
$ws.Range("C9").Value = "heavy"
$ws.Range("D9").Formula = "=(F9*3600)/E9"
$ws.Range("E9").Value = 150
$ws.Range("F9").Value = 1.08
$ws.Range("G9").Value = 0.998
$ws.Range("H9").Value = 0.999
$ws.Range("I9").Value = 0.995
$ws.Range("J9").Value = 0.961

# Row 10
$ws.Range("A10").Value = "11n"
$ws.Range("B10").Value = 256
$ws.Range("C10").Value = "heavy"
$ws.Range("D10").Formula = "=(F10*3600)/E10"
$ws.Range("E10").Value = 150
$ws.Range("F10").Value = 0.477
$ws.Range("G10").Value = 0.993
$ws.Range("H10").Value = 0.992
$ws.Range("I10").Value = 0.994
$ws.Range("J10").Value = 0.918

# Row 11
$ws.Range("A11").Value = "11m"
$ws.Range("B11").Value = 640
$ws.Range("C11").Value = "no"
$ws.Range("D11").Formula = "=(F11*3600)/E11"
$ws.Range("E11").Value = 52
$ws.Range("F11").Value = 0.7
$ws.Range("G11").Value = 0.99
$ws.Range("H11").Value = 0.994
$ws.Range("I11").Value = 0.994
$ws.Range("J11").Value = 0.953
